# Update countries & provincias Spain
# Daily COVID data refresh: Ucrania and Oman overtook neighbouring rows in
# the "Casos totales" ranking, so their rows moved up (pushing the rows
# below them down, since the sheet is sorted by total cases). The country
# label cells (column A) keep referencing the same shared-string slot, but
# since that slot's position in the ranking changed, the text now shown in
# that row changes too; the numeric columns (B:H) simply carry the data
# for whichever country now sits at that row. Estado de Palestina (row 108)
# also received a small update, and the "last updated" footer timestamp
# moved from 08:22 to 08:52.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Footer timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 08:52"

# --- Ucrania overtakes Colombia (rows 50-51) ---------------------------
# Row 50 becomes Ucrania's (new, higher) numbers.
$ws.Range("A50").Value = "Ucrania"
$ws.Range("B50").Value = 2511
$ws.Range("C50").Value = 308
$ws.Range("D50").Value = 79
$ws.Range("E50").Value = 2359
$ws.Range("F50").Value = 33
$ws.Range("G50").Value = 4
$ws.Range("H50").Value = 73

# Row 51 becomes Colombia's (unchanged) numbers.
$ws.Range("A51").Value = "Colombia"
$ws.Range("B51").Value = 2473
$ws.Range("C51").Value = 0
$ws.Range("D51").Value = 197
$ws.Range("E51").Value = 2196
$ws.Range("F51").Value = 85
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 80

# --- Oman overtakes Afganistan and Uruguay (rows 89-91) -----------------
# Row 89 becomes Oman's (new, higher) numbers.
$ws.Range("A89").Value = "Oman"
$ws.Range("B89").Value = 546
$ws.Range("C89").Value = 62
$ws.Range("D89").Value = 109
$ws.Range("E89").Value = 434
$ws.Range("F89").Value = 3
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 3

# Row 90 becomes Afganistan's (unchanged) numbers.
$ws.Range("A90").Value = "Afganistan"
$ws.Range("B90").Value = 521
$ws.Range("C90").Value = 0
$ws.Range("D90").Value = 32
$ws.Range("E90").Value = 474
$ws.Range("F90").Value = 0
$ws.Range("G90").Value = 0
$ws.Range("H90").Value = 15

# Row 91 becomes Uruguay's (unchanged) numbers.
$ws.Range("A91").Value = "Uruguay"
$ws.Range("B91").Value = 494
$ws.Range("C91").Value = 21
$ws.Range("D91").Value = 214
$ws.Range("E91").Value = 273
$ws.Range("F91").Value = 15
$ws.Range("G91").Value = 0
$ws.Range("H91").Value = 7

# --- Estado de Palestina small update (row 108) --------------------------
$ws.Range("B108").Value = 268
$ws.Range("C108").Value = 1
$ws.Range("D108").Value = 46
